$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Abril de 2020 a las 19:22"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 748113
$ws.Cells.Item(4, 3).Value = 9321
$ws.Cells.Item(4, 4).Value = 68822
$ws.Cells.Item(4, 5).Value = 639442
$ws.Cells.Item(4, 6).Value = 13551
$ws.Cells.Item(4, 7).Value = 835
$ws.Cells.Item(4, 8).Value = 39849

# Row 8
$ws.Cells.Item(8, 1).Value = "Alemania"
$ws.Cells.Item(8, 2).Value = 144387
$ws.Cells.Item(8, 3).Value = 663
$ws.Cells.Item(8, 4).Value = 88000
$ws.Cells.Item(8, 5).Value = 51840
$ws.Cells.Item(8, 6).Value = 2889
$ws.Cells.Item(8, 7).Value = 9
$ws.Cells.Item(8, 8).Value = 4547

# Row 18
$ws.Cells.Item(18, 1).Value = "Suiza"
$ws.Cells.Item(18, 2).Value = 27740
$ws.Cells.Item(18, 3).Value = 336
$ws.Cells.Item(18, 4).Value = 17800
$ws.Cells.Item(18, 5).Value = 8547
$ws.Cells.Item(18, 6).Value = 386
$ws.Cells.Item(18, 7).Value = 25
$ws.Cells.Item(18, 8).Value = 1393

# Row 21
$ws.Cells.Item(21, 1).Value = "Irlanda"
$ws.Cells.Item(21, 2).Value = 15251
$ws.Cells.Item(21, 3).Value = 493
$ws.Cells.Item(21, 4).Value = 77
$ws.Cells.Item(21, 5).Value = 14564
$ws.Cells.Item(21, 6).Value = 294
$ws.Cells.Item(21, 7).Value = 39
$ws.Cells.Item(21, 8).Value = 610

# Row 69
$ws.Cells.Item(69, 1).Value = "Irak"
$ws.Cells.Item(69, 2).Value = 1539
$ws.Cells.Item(69, 3).Value = 26
$ws.Cells.Item(69, 4).Value = 1009
$ws.Cells.Item(69, 5).Value = 448
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 82

# Row 70
$ws.Cells.Item(70, 1).Value = "Estonia"
$ws.Cells.Item(70, 2).Value = 1528
$ws.Cells.Item(70, 3).Value = 16
$ws.Cells.Item(70, 4).Value = 164
$ws.Cells.Item(70, 5).Value = 1324
$ws.Cells.Item(70, 6).Value = 10
$ws.Cells.Item(70, 7).Value = 2
$ws.Cells.Item(70, 8).Value = 40

# Row 86
$ws.Cells.Item(86, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(86, 2).Value = 846
$ws.Cells.Item(86, 3).Value = 114
$ws.Cells.Item(86, 4).Value = 102
$ws.Cells.Item(86, 5).Value = 742
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 2

# Row 87
$ws.Cells.Item(87, 1).Value = "Ghana"
$ws.Cells.Item(87, 2).Value = 834
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 99
$ws.Cells.Item(87, 5).Value = 726
$ws.Cells.Item(87, 6).Value = 4
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 9

# Row 88
$ws.Cells.Item(88, 1).Value = "Costa de Marfil"
$ws.Cells.Item(88, 2).Value = 801
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 239
$ws.Cells.Item(88, 5).Value = 554
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 8

# Row 89
$ws.Cells.Item(89, 1).Value = "Republica de Chipre"
$ws.Cells.Item(89, 2).Value = 767
$ws.Cells.Item(89, 3).Value = 6
$ws.Cells.Item(89, 4).Value = 81
$ws.Cells.Item(89, 5).Value = 674
$ws.Cells.Item(89, 6).Value = 15
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 12

# Row 93
$ws.Cells.Item(93, 1).Value = "Libano"
$ws.Cells.Item(93, 2).Value = 673
$ws.Cells.Item(93, 3).Value = 1
$ws.Cells.Item(93, 4).Value = 102
$ws.Cells.Item(93, 5).Value = 550
$ws.Cells.Item(93, 6).Value = 27
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 21

# Row 144
$ws.Cells.Item(144, 1).Value = "Liberia"
$ws.Cells.Item(144, 2).Value = 91
$ws.Cells.Item(144, 3).Value = 15
$ws.Cells.Item(144, 4).Value = 7
$ws.Cells.Item(144, 5).Value = 76
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 8

# Row 145
$ws.Cells.Item(145, 1).Value = "Bermudas"
$ws.Cells.Item(145, 2).Value = 86
$ws.Cells.Item(145, 3).Value = 3
$ws.Cells.Item(145, 4).Value = 35
$ws.Cells.Item(145, 5).Value = 46
$ws.Cells.Item(145, 6).Value = 10
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 5

# Row 146
$ws.Cells.Item(146, 1).Value = "Togo"
$ws.Cells.Item(146, 2).Value = 84
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 52
$ws.Cells.Item(146, 5).Value = 27
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 5

# Row 147
$ws.Cells.Item(147, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(147, 2).Value = 79
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 4
$ws.Cells.Item(147, 5).Value = 75
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 0

# Row 148
$ws.Cells.Item(148, 1).Value = "Liechtenstein"
$ws.Cells.Item(148, 2).Value = 79
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 55
$ws.Cells.Item(148, 5).Value = 23
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 1
